$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Math book"
$ws.Range("B1").Value = 12525
$ws.Range("C1").Value = 4

$ws.Range("A2").Value = "AI Book"
$ws.Range("B2").Value = 3589
$ws.Range("C2").Value = 6

$ws.Range("C2").Select()
